# Bump the ObjTables/SBtab schema metadata embedded in each sheet's header row(s):
#   objTablesVersion '0.0.9' -> '1.0.0', date '2020-04-27 01:10:15/16' -> '2020-05-29 00:24:11'
$wb = $excel.ActiveWorkbook
$newDate = "2020-05-29 00:24:11"
$newVer = "1.0.0"

# All sheets are protected (no password); unprotect so the header cells can be edited
foreach ($sheet in $wb.Worksheets) {
    $sheet.Unprotect()
}

# Row 1 of the first sheet carries the workbook-wide "!!!ObjTables" banner
$wsCompartment = $wb.Worksheets.Item("!!Compartment")
$wsCompartment.Range("A1").Value = "!!!ObjTables schema='SBtab' objTablesVersion='$newVer' date='$newDate'"

# Every sheet (including !!Compartment, on row 2) carries its own "!!ObjTables ... class='<Name>'" header
$classNames = @(
    "Compartment",
    "Compound",
    "Definition",
    "Enzyme",
    "FbcObjective",
    "Gene",
    "Layout",
    "Measurement",
    "PbConfig",
    "Position",
    "Protein",
    "Quantity",
    "QuantityInfo",
    "QuantityMatrix",
    "Reaction",
    "ReactionStoichiometry",
    "Regulator",
    "Relation",
    "Relationship",
    "SparseMatrix",
    "SparseMatrixColumn",
    "SparseMatrixOrdered",
    "SparseMatrixRow",
    "StoichiometricMatrix",
    "rxnconContingencyList",
    "rxnconReactionList"
)

foreach ($className in $classNames) {
    $ws = $wb.Worksheets.Item("!!" + $className)
    $headerRow = 1
    if ($className -eq "Compartment") { $headerRow = 2 }
    $header = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='" + $className + "' name='" + $className + "' date='$newDate' objTablesVersion='$newVer'"
    $ws.Cells.Item($headerRow, 1).Value = $header
}

# Restore sheet protection to match the pre-edit state
foreach ($sheet in $wb.Worksheets) {
    $sheet.Protect()
}
